$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 618
$ws1.Range("F13").Value = 1230
$ws1.Range("F14").Value = 20
$ws1.Range("F15").Value = 2857
$ws1.Range("F16").Value = 442
$ws1.Range("F17").Value = 543

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 618
$ws4.Range("F14").Value = 1230
$ws4.Range("F15").Value = 20
$ws4.Range("F16").Value = 2857
$ws4.Range("F17").Value = 442
$ws4.Range("F18").Value = 543
